$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '57.361.07'
$c.Style = 'Normal'
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  -0.89%  '
$c.Style = 'Normal'

$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '3.108.29'
$c.Style = 'Normal'
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  +0.12%  '
$c.Style = 'Normal'

$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c.Style = 'Normal'

$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '525.71'
$c.Style = 'Normal'
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  -0.06%  '
$c.Style = 'Normal'

$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '138.10'
$c.Style = 'Normal'
$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  -2.62%  '
$c.Style = 'Normal'

$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '
$c.Style = 'Normal'

$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '3.106.38'
$c.Style = 'Normal'
$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  +0.08%  '
$c.Style = 'Normal'

$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  +2.39%  '
$c.Style = 'Normal'

$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '7.36'
$c.Style = 'Normal'
$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  +2.04%  '
$c.Style = 'Normal'

$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  -0.96%  '
$c.Style = 'Normal'

$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  +2.57%  '
$c.Style = 'Normal'

$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '3.640.33'
$c.Style = 'Normal'
$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  -0.03%  '
$c.Style = 'Normal'

$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  +1.53%  '
$c.Style = 'Normal'

$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '25.59'
$c.Style = 'Normal'
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  +0.08%  '
$c.Style = 'Normal'

$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  -0.94%  '
$c.Style = 'Normal'

$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '57.492.54'
$c.Style = 'Normal'

$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '3.108.79'
$c.Style = 'Normal'
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  +0.22%  '
$c.Style = 'Normal'

$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '5.93'
$c.Style = 'Normal'
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  -2.87%  '
$c.Style = 'Normal'

$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  -1.68%  '
$c.Style = 'Normal'

$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '7.94'
$c.Style = 'Normal'
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  -0.74%  '
$c.Style = 'Normal'

$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '350.90'
$c.Style = 'Normal'
$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  +2.54%  '
$c.Style = 'Normal'

$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  +0.14%  '
$c.Style = 'Normal'

$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '68.42'
$c.Style = 'Normal'
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  +1.69%  '
$c.Style = 'Normal'

$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  -1.89%  '
$c.Style = 'Normal'

$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  -0.57%  '
$c.Style = 'Normal'

$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  -0.18%  '
$c.Style = 'Normal'

$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '0.0₃0893'
$c.Style = 'Normal'
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  -2.79%  '
$c.Style = 'Normal'

$c = $ws.Range("B29")
$c.NumberFormat = '@'
$c.Value = 'InternetComputer(DFINITY)'
$c.Style = 'Normal'
$c = $ws.Range("C29")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c.Style = 'Normal'
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '7.40'
$c.Style = 'Normal'
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  +1.52%  '
$c.Style = 'Normal'

$c = $ws.Range("B30")
$c.NumberFormat = '@'
$c.Value = 'USDe'
$c.Style = 'Normal'
$c = $ws.Range("C30")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c.Style = 'Normal'
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.Style = 'Normal'

$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  +0.40%  '
$c.Style = 'Normal'

$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '6.00'
$c.Style = 'Normal'
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  -7.01%  '
$c.Style = 'Normal'

$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '20.89'
$c.Style = 'Normal'
$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  -0.51%  '
$c.Style = 'Normal'

$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '5.04'
$c.Style = 'Normal'
$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  +8.42%  '
$c.Style = 'Normal'

$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  -3.47%  '
$c.Style = 'Normal'

$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '159.72'
$c.Style = 'Normal'
$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  +0.95%  '
$c.Style = 'Normal'

$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  -1.53%  '
$c.Style = 'Normal'

$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '26.44'
$c.Style = 'Normal'
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  +0.41%  '
$c.Style = 'Normal'

$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  -0.17%  '
$c.Style = 'Normal'

$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '0.0659'
$c.Style = 'Normal'
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  -1.15%  '
$c.Style = 'Normal'

$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '1.60'
$c.Style = 'Normal'
$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  +0.50%  '
$c.Style = 'Normal'

$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  +1.67%  '
$c.Style = 'Normal'

$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  +2.06%  '
$c.Style = 'Normal'

$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '2.409.09'
$c.Style = 'Normal'
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  +5.46%  '
$c.Style = 'Normal'

$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '36.70'
$c.Style = 'Normal'
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  -0.45%  '
$c.Style = 'Normal'

$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  -0.05%  '
$c.Style = 'Normal'

$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '3.148.48'
$c.Style = 'Normal'
$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.Style = 'Normal'

$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  +0.65%  '
$c.Style = 'Normal'

$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  -2.24%  '
$c.Style = 'Normal'

$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  -1.79%  '
$c.Style = 'Normal'

$c = $ws.Range("B51")
$c.NumberFormat = '@'
$c.Value = 'SuiNetwork'
$c.Style = 'Normal'
$c = $ws.Range("C51")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c.Style = 'Normal'
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '0.769'
$c.Style = 'Normal'
$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  +3.10%  '
$c.Style = 'Normal'
